$wb = $excel.ActiveWorkbook

# --- "network_optimized_weights" sheet: was the active tab (tabSelected),
# move the selection there to H6, it will lose tabSelected once another
# sheet becomes the active tab below.
$wsWeights = $wb.Worksheets.Item("network_optimized_weights")
$wsWeights.Select()
$wsWeights.Range("H6").Select()

# --- "optimization_parameters" sheet: insert a new row 9 ("L_curve"),
# rename the old "Model" row label to "production_function", and make
# this sheet the active tab with selection A9:B9.
$ws = $wb.Worksheets.Item("optimization_parameters")
$ws.Select()

$ws.Rows.Item(9).Insert()

$ws.Range("A8").Value = "production_function"
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 0

$ws.Range("A9:B9").Select()
